$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# A26 label text: "run number" -> "run_num"
$ws.Range("A26").Value = "run_num"

# B25 value: 1 -> 0
$ws.Range("B25").Value = 0

# Update the view's active selection (was B33, now K34) and scroll
# position (top-left visible cell moves from A23 up to A22).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("K34").Select()
